# Weekly crime-stat refresh: CompStat 103rd Precinct report
# Updates report week (Volume/date header) and the weekly/28-day/YTD/
# 2-year crime-complaint figures (and their computed % changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/issue number and report week date range ---
$ws.Range("A8").Value = "Volume 33   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/2/2026  Through  2/8/2026"

# --- Helper: turn a cell into the "0" / "***.*" text placeholder used
#     throughout this sheet when a count/percentage is not applicable,
#     copying both the display format and the exact text from a cell
#     elsewhere on the sheet that already holds it. ---
function Set-PlaceholderText($destRef, $donorRef) {
    $ws.Range($donorRef).Copy()
    $ws.Range($destRef).PasteSpecial(-4122)
    $ws.Range($donorRef).Copy()
    $ws.Paste($ws.Range($destRef))
}

# --- Helper: turn a placeholder-text cell back into a real number,
#     copying the numeric/percentage display format from a stable
#     donor cell of the same column-kind before writing the value. ---
function Set-NumberFromText($destRef, $donorRef, $value) {
    $ws.Range($donorRef).Copy()
    $ws.Range($destRef).PasteSpecial(-4122)
    $ws.Range($destRef).Value = $value
}

# --- Cells that become "no data" placeholders this week ---
Set-PlaceholderText "D14" "C14"
Set-PlaceholderText "E14" "L14"
Set-PlaceholderText "D15" "C14"
Set-PlaceholderText "E15" "L14"
Set-PlaceholderText "C22" "C14"
Set-PlaceholderText "D27" "C14"
Set-PlaceholderText "E27" "L14"

# --- Cells that were placeholders last week but have real figures now ---
Set-NumberFromText "D22" "C15" 1
Set-NumberFromText "E22" "N14" -100
Set-NumberFromText "G22" "C15" 1
Set-NumberFromText "H22" "N14" 300
Set-NumberFromText "J22" "C15" 1
Set-NumberFromText "K22" "N14" 500
Set-NumberFromText "M29" "N14" -100
Set-NumberFromText "M30" "N14" -100

# --- Updated weekly/28-day/YTD/2-year crime counts and % changes ---
# Row 15
$ws.Range("F15").Value = 5
$ws.Range("H15").Value = 66.666666666666
$ws.Range("I15").Value = 8
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = 300
$ws.Range("M15").Value = 166.666666666667
$ws.Range("N15").Value = 14.285714285714
# Row 16
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("G16").Value = 24
$ws.Range("H16").Value = -20.833333333333
$ws.Range("I16").Value = 27
$ws.Range("J16").Value = 32
$ws.Range("K16").Value = -15.625
$ws.Range("L16").Value = -12.903225806451
$ws.Range("M16").Value = -43.75
$ws.Range("N16").Value = -85.789473684210
# Row 17
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 50
$ws.Range("G17").Value = 45
$ws.Range("H17").Value = 11.111111111111
$ws.Range("I17").Value = 71
$ws.Range("J17").Value = 74
$ws.Range("K17").Value = -4.054054054054
$ws.Range("L17").Value = 16.393442622950
$ws.Range("M17").Value = 195.833333333333
$ws.Range("N17").Value = -15.476190476190
# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -15.384615384615
$ws.Range("I18").Value = 22
$ws.Range("J18").Value = 22
$ws.Range("L18").Value = 57.142857142857
$ws.Range("M18").Value = -8.333333333333
$ws.Range("N18").Value = -82.8125
# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -27.272727272727
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = 13.513513513513
$ws.Range("I19").Value = 55
$ws.Range("J19").Value = 61
$ws.Range("K19").Value = -9.836065573770
$ws.Range("L19").Value = -19.117647058823
$ws.Range("M19").Value = 17.021276595744
$ws.Range("N19").Value = -51.327433628318
# Row 20
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 21
$ws.Range("J20").Value = 19
$ws.Range("K20").Value = 10.526315789473
$ws.Range("L20").Value = -30
$ws.Range("M20").Value = 40
$ws.Range("N20").Value = -85.517241379310
# Row 21
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = -2.941176470588
$ws.Range("F21").Value = 143
$ws.Range("G21").Value = 139
$ws.Range("H21").Value = 2.877697841726
$ws.Range("I21").Value = 204
$ws.Range("J21").Value = 213
$ws.Range("K21").Value = -4.225352112676
$ws.Range("L21").Value = -0.970873786407
$ws.Range("M21").Value = 26.708074534161
$ws.Range("N21").Value = -69.461077844311
# Row 22
$ws.Range("F22").Value = 4
$ws.Range("M22").Value = 100
# Row 23
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -66.666666666666
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 8
$ws.Range("K23").Value = -50
$ws.Range("L23").Value = -42.857142857142
$ws.Range("M23").Value = -20
# Row 24
$ws.Range("C24").Value = 47
$ws.Range("D24").Value = 39
$ws.Range("E24").Value = 20.512820512820
$ws.Range("F24").Value = 154
$ws.Range("G24").Value = 135
$ws.Range("H24").Value = 14.074074074074
$ws.Range("I24").Value = 217
$ws.Range("J24").Value = 198
$ws.Range("K24").Value = 9.595959595959
$ws.Range("L24").Value = 24.712643678160
$ws.Range("M24").Value = 99.082568807339
# Row 25
$ws.Range("C25").Value = 30
$ws.Range("D25").Value = 29
$ws.Range("E25").Value = 3.448275862068
$ws.Range("F25").Value = 98
$ws.Range("G25").Value = 76
$ws.Range("H25").Value = 28.947368421052
$ws.Range("I25").Value = 128
$ws.Range("J25").Value = 110
$ws.Range("K25").Value = 16.363636363636
$ws.Range("L25").Value = 34.736842105263
# Row 26
$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 17
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 75
$ws.Range("G26").Value = 72
$ws.Range("H26").Value = 4.166666666666
$ws.Range("I26").Value = 105
$ws.Range("J26").Value = 93
$ws.Range("K26").Value = 12.903225806451
$ws.Range("L26").Value = 9.375
$ws.Range("M26").Value = 9.375
# Row 27
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 8
$ws.Range("H27").Value = 166.666666666667
$ws.Range("I27").Value = 12
$ws.Range("K27").Value = 200
$ws.Range("L27").Value = 140
# Row 28
$ws.Range("C28").Value = 6
$ws.Range("E28").Value = 500
$ws.Range("F28").Value = 12
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 13
$ws.Range("J28").Value = 8
$ws.Range("K28").Value = 62.5
$ws.Range("L28").Value = 30
# Row 31
$ws.Range("F31").Value = 3
$ws.Range("I31").Value = 3
